$wb = $excel.ActiveWorkbook

# --- YDS sheet: extend R/P OFF/DEF logs with Week 16 data ---
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = "11 4 1 3 10 6 3 2 0 11 2 3 3 2 5 3 1 0 2 1 5 2 8 4 5 2 0 4 0 0 5 13 4 5 1 4 3 11 1 -7 6 1 0 6 7 14 4 3 7 4 2 2 14 4 3 6 -3 1 8 -2 9 3 2 9 2 0 5 3 10 0 0 7 4 2 4 0 9 1 -3 16 2 -3 0 7 12 5 6 5 2 1 3 1 5 12 5 3 5 5 0 2 3 0 2 1 1 9 4 5 3 1 12 4 6 10 5 4 2 0 0 6 2 23 6 -3 7 5 3 9 3 4 2 -1 18 3 3 0 34 0 2 0 2 -3 5 14 5 2 9 4 2 3 5 20 4 10 94 2 17 3 5 1 14 6 3 -2 2 0 4 -4 9 1 2 0 2 3 9 3 1 4 17 6 0 1 0 15 4 0 6 5 6 2 9 13 8 11 3 6 4 3 3 16 4 8 21 15 2 9 5 7 8 17 8 3 1 5 2 5 0 3 0 6 2 2 4 0 3 5 3 1 0 5 6 1 26 0 0 2 2 1 1 5 3 5 -4 8 0 20 2 8 16 -1 -1 1 5 5 9 0 1 3 9 12 6 4 11 16 3 7 2 8 2 -1 7 0 4 -5 5 2 4 8 8 3 0 21 3 3 3 3 0 2 5 1 0 2 -2 24 1 11 3 1 12 4 3 8 4 4 29 6 3 9 12 2 12 6 2 10 6 10 5 1 1 31 8 5 12 5 4 5 11 3 3 1 -1 1 14 0 2 12 7 4 6 2 1 2 2 3 1 3 4 1 3 7 2 -3 5 0 10 3 -1 9 3 9 5 0 7 1 1 2 2 6 5 6 2 6 6 5 7 8 22 5 4 36 1 47 3 10 6 1 2 5 18 2 3 1 20 2 0 4 -2 3 6 18 2 -4 -1 2 -4 4 5 9 5 3 12 8 13 6 7 17 2 8 10 5 1 -2 1 5 4 5 33 3 3 4 -1 3 3 4 8 9 4 -1 0 5 9 4 3 1 5 1 5 3 10 -2 3 4 5 45 3 4 3 10 3 6 8 1 9 5 3 3 8 -3 2 2 8 52 10 16 28 -3 4 14 -2 0 6 8 9 4 45 0 5 3 1 6 6 -1 1 3 5 6 1 4 11 7 1 3 5 4 4 2 6 1 1 -3 5 4 1 1 0 4 2 4 1 3 2 0 8 3 0 -3 -1 2 7 -1 2 2 0 0 1 16 1 19 4 8 8 -1 4 7 8 1 2 1 -2 2 7 0 1 4 7 3 2 4 8 4 1 11 3 4 -2 9 2 6 7 0 60 3 6 2 19 0 3 8 4 4 1 0 12 4 1 4 1 3 4 17 0 19 1 -1 4 2 1 3 4 8 2 8 2 9 5 1 2 3 2 28 10 9 2 2 12 1 5 1 -1 13 0 7 4 0 2 2 2 4 16 1 7 2 0 4 3 5 2 3 1 2 11 2 3 18 9 11 22 1 2 2 7 5 5 1 0 -4 2 4 4 4 6 -1 0 1 5 9 5 11 1 6 9 4 2 3 3 8 7 13 9 1 4 -2 0 4 15 0 9 5 4 11 10 -1 2 4 3 76 2 -3 7 4 6 3 19 0 4 3 1 2 3 0 -1 13 3 9 2 3 4 0 8 5 0 3 4 -1 2 4 0 11 1 -1 2 2 4 6 2 1 2 8 5 4 -4 9 3 2 -4 5 3 1 3 4 1 13 2 5 9 2 1 0 6 1 1 2 2 -1 4 2 4 5 3 3 6 -1 -3 2 3 1 3 -3 8 -7 0 1 4 2 14 2 11 6 3 4 2 1 6 0 3 2 1 1 2 4 2 5 4 -1 6 5 7 1 2 1 5 3 2 2 1 1 0 2 3 4 5 6 -5 0 1 11 7 0 16 4 12 5 4 9 -1 8 -1 9 -1 1 2 16 2 2 1 -3 -1 4 5 2 11 2 5 2 1 4 2 17 3 13 3 68 6 30 5 6 3 7 14 4 2 8 4 3 1 2 1 5 12 9 2 3 1 2 6 3 17 0 5 4 4 6 8 2 5 3 4 3 2 -5 3 1 1 -1 2 4 5 6 1 1 12 0 2 -1 3 6 2 1 3 8 9 2 1 2 3 6 0 3 9 1 2 0 9 12 1 6 3 10 7 2 4 4 15 4 5 2 11 0 20 6 7 3 1 4 3 7 0 0 8 9 -1 0 0 3 1 0 5 10 3 2 4 3 -1 4 5 23 5 1 -1"
$ydsWs.Range("B3").Value = "22 16 20 4 5 1 23 10 13 9 12 4 7 5 7 6 -1 8 7 4 13 7 5 5 8 9 6 11 6 63 13 14 9 8 12 10 10 8 5 4 16 17 4 18 18 5 5 44 -1 16 0 13 10 9 9 12 23 5 6 5 5 21 6 38 61 6 8 11 14 8 1 16 11 20 16 14 7 6 9 12 6 3 3 12 15 6 6 6 4 4 12 7 6 18 13 10 7 9 10 8 6 2 13 5 20 22 45 12 7 15 -1 8 11 5 6 23 7 8 9 7 53 2 13 19 7 4 4 3 73 4 6 20 9 17 9 -1 7 -2 21 7 10 5 27 24 9 23 19 0 6 9 11 9 10 26 12 10 12 38 6 17 9 40 6 12 30 2 -2 13 21 20 5 8 5 3 6 18 6 5 7 4 8 18 14 7 8 1 -1 7 6 2 12 50 5 12 16 9 5 9 22 25 14 14 11 11 12 4 7 8 3 69 12 21 16 10 12 10 18 12 40 17 1 3 13 14 19 43 18 22 18 9 11 5 6 4 14 6 11 8 3 14 9 13 37 11 25 6 6 17 13 5 10 5 19 9 13 9 5 17 16 0 6 6 19 10 10 75 7 14 14 16 5 5 5 -1 20 14 7 24 2 12 3 12 10 12 12 8 12 9 11 10 19 6 21 6 8 1 19 4 4 29 12 8 6 11 11 5 5 5 9 52 28 4 14 5 9 10 35 5 10 6 7 18 12 3 -2 10 6 5 3 4 8 12 11 39 10 9 9 4 4 19 13 9 10 5 10 8 10 8 7 10 6 14 13 5 17 51 15 15 19 6 15 15 22 12 12 5 15 16 10 14 4 9 12 4 5 6 6 6 9 3 25 7 18 26 14 13 13 -2 10 11 13 7 3 27 5 11 14 16 4 7 10 0 12 7 10 9 1 23 9 14 18 9 15 12 11 6 6 2 20 -6 5 9 12 10 12 25 10 22 22 8 14 26 20 6 7 1 14 11 5 48 7 8 8 8 4 13 20 6 6 14 22 5 8 9 10 18 3 5 3 24 8 13 22 10 2 20 6 12 9 11 6 5 8 24 8 7 3 8 8 6 14 7 57 13 10 11 7 5 10 11 14 9 10 10 13 11 13 15 5 -3 7 2 10 14 1 8 3 5 9 9 4 7 11 7 13 16 3 2 50 9 16 17 10 1 16 16 9 10 -1 2 5 7 2 39 0 11 10 10 11 15 1 8 7 12 7 10 8 12 5 11 -1 16 9 5 7 18 -2 4 7 46 11 5 9 9 4 4 8 14 7 8 10 2 14 8 1 2 3 24 7 14 9 6 11 15 1 16 9 12 5 4 8 11 4 10 24 9 5 13 7 12 7 5 11 12 1 4 8 5 18 7 9 13 11 5 -4 -3 3 5 12 5 9 8 -8 14 3 7 -3 11 4 0 12 0 24 18 -1 13 2 42 12 18 18 13 6 4"
$ydsWs.Range("C2").Value = "7 0 4 7 0 3 10 3 6 3 3 1 9 7 1 0 -2 0 25 3 0 1 12 7 0 2 3 3 5 2 0 0 0 2 4 14 3 2 3 1 4 5 6 2 4 39 12 9 5 3 1 2 11 2 9 3 11 4 2 39 3 1 4 4 1 3 3 -2 12 5 4 39 -3 2 0 13 1 14 6 1 10 15 1 7 12 0 3 3 3 -3 -2 7 6 3 3 4 1 2 4 2 5 7 2 2 6 4 34 1 4 7 -1 4 -2 10 6 4 11 0 11 1 4 3 4 1 -2 2 -2 3 3 11 8 0 0 4 2 6 1 6 10 7 3 9 2 1 24 8 2 9 -1 6 -3 0 5 0 3 5 3 6 9 2 6 2 1 1 12 -3 4 0 12 4 7 2 4 5 2 1 16 1 -1 6 3 1 2 0 3 1 3 5 1 1 2 -1 5 2 2 6 0 -2 -3 4 3 3 5 5 1 4 4 0 3 3 9 12 3 2 8 4 3 14 13 -3 0 0 4 21 2 3 2 2 18 4 -2 3 2 5 0 2 10 7 1 2 6 9 -3 1 11 4 7 2 2 1 3 2 3 2 3 2 5 0 1 8 8 14 2 2 5 2 2 3 3 3 10 1 1 3 11 2 2 2 0 2 1 0 3 14 8 0 4 -1 -1 -1 0 2 3 4 3 4 4 1 1 2 1 4 29 0 2 1 2 8 5 2 4 0 0 5 3 1 1 -1 3 5 -4 1 2 2 3 3 4 1 3 9 47 13 7 3 7 8 3 -1 3 4 5 2 12 4 0 1 2 7 11 4 7 2 0 10 -1 4 11 6 3 12 5 1 14 1 6 9 14 1 6 2 5 6 -3 5 59 -1 8 6 -3 8 1 0 10 0 6 8 9 0 0 2 7 4 3 1 1 2 3 1 6 10 5 3 1 19 26 1 7 3 3 2 4 0 3 3 5 0 4 8 1 15 2 48 8 10 -1 9 23 2 4 4 -2 1 3 12 2 14 4 33 5 4 5 10 0 4 3 10 7 11 8 4 8 0 0 6 2 4 3 7 16 4 2 6 1 5 3 2 1 2 6 8 2 1 1 1 2 7 5 1 0 2 4 4 2 25 1 8 1 -3 -1 5 2 12 2 23 -1 2 2 6 9 -7 5 6 6 1 5 9 5 2 2 1 4 1 2 -2 9 6 2 5 -1 5 -2 -2 2 2 6 13 0 2 8 5 2 -3 0 5 58 7 4 1 3 1 6 -1 9 0 5 10 12 2 9 2 10 8 9 8 5 2 1 4 3 4 8 6 1 3 14 3 7 3 -3 7 5 2 1 7 5 2 3 4 0 4 5 0 4 5 5 2 2 11 -1 13 18 6 9 9 6 2 7 3 12 11 0 32 2 6 5 4 2 3 1 6 5 4 1 4 4 5 3 4 1 14 4 12 2 5 2 11 4 6 0 4 3 0 8 0 1 2 5 10 3 1 6 8 0 2 2 6 4 4 6 13 3 0 5 3 2 5 2 7 4 4 2 1 -3 2 3 2 3 7 -3 4 5 8 1 3 7 1 1 -1 2 0 5 -2 3 1 9 0 1 -1 1 0 -1 0 3 1 9 0 2 2 2 2 2 1 6 5 7 4 2 -1 4 6 10 19 -1 7 14 -4 3 5 2 2 -1 0 1 2 1 10 3 0 7 2 7 0 0 1 2 3 -3 6 -5 4 8 3 1 9 2 0 1 7 5 3 1 0 7 3 4 0 5 4 12 4"
$ydsWs.Range("C3").Value = "17 3 12 8 -6 9 6 2 16 31 11 8 25 7 -2 21 -2 15 12 5 18 12 11 9 45 19 17 12 7 8 14 -3 9 7 23 6 4 8 13 9 7 10 6 12 6 2 24 8 14 11 9 11 9 12 8 6 3 31 16 33 15 71 7 11 3 7 8 9 19 15 16 9 3 15 14 12 12 8 5 9 6 10 7 -7 14 6 11 10 6 8 14 10 22 9 2 11 22 7 18 17 6 16 8 17 6 12 6 11 35 0 2 4 19 5 6 4 22 53 15 7 6 9 7 14 6 9 4 5 14 9 5 11 22 2 -5 21 18 6 9 7 7 5 2 7 9 2 7 5 7 28 7 4 3 8 9 0 10 24 2 6 6 7 8 13 2 13 3 7 12 8 12 13 9 13 7 4 18 6 15 8 24 4 4 5 4 14 9 7 6 34 5 5 7 22 18 3 -4 23 4 5 5 8 5 27 13 1 6 4 16 1 2 12 11 8 11 5 12 6 9 9 18 7 11 30 14 13 11 3 21 8 8 11 2 14 11 2 7 3 40 6 6 5 14 30 5 4 6 6 31 14 5 3 6 8 15 16 9 9 12 21 12 11 19 19 14 5 10 4 1 12 14 9 11 50 23 14 7 8 7 3 4 5 10 35 5 4 26 6 7 8 7 6 1 75 11 21 5 12 17 17 13 8 8 16 7 13 0 9 9 4 11 4 22 4 9 6 2 5 5 -5 2 10 18 5 1 13 4 11 22 26 27 12 12 12 6 -2 8 0 6 4 13 11 12 2 6 5 8 2 5 44 7 6 39 12 14 7 36 2 11 7 21 11 5 11 7 8 13 7 5 -1 7 17 21 6 7 23 5 9 6 8 11 8 5 27 32 20 4 8 42 9 2 14 7 7 12 8 14 6 7 2 11 14 38 20 27 0 4 9 7 7 32 5 29 7 17 5 28 9 19 -6 11 4 6 9 3 17 20 9 17 4 38 1 21 8 5 3 10 3 17 18 4 22 3 12 36 3 0 10 26 11 9 51 4 9 4 5 16 63 11 8 12 8 1 14 10 11 1 7 68 3 7 3 27 10 8 8 7 12 4 18 4 3 11 4 5 27 36 18 3 4 9 3 9 -4 13 30 8 6 3 6 15 6 54 5 29 3 53 4 0 9 29 11 8 5 18 5 -2 25 7 22 13 7 13 7 18 58 9 4 8 4 5 12 19 7 8 9 12 12 6 4 3 11 6 7 3 3 21 19 2 14 5 17 12 29 1 5 9 11 18 9 15 3 4 1 6 6 31 12 7 20 7 8 6 6 5 10 9 23 8 11 24 12 14 6 4 9 10 4 2 37 5 8 0 7 10 7 7 7 8 9 7 12 5 6 1 9 6 5 -2 11 3 9 38 8 8 14 12 4 0 7 1 10 6 18 20 26 15 6 17 6 3 4 9 22 13 7 2 13 3 16 5 9 6 9 3 10 8 6 7 6 46 11 3 12 21 21 8 10 34 26 22 5 17 15 8 10 9 8 12 3 13 10 5 6 2 8 7 6 11 4 22 20 16 13 12 38 5 19 5 3 21 10 7 20 9 11 5 41 4 7 7 8 14 3 7 5 7 3 12 3 31 7 14 6 5 5 9 9 31 41 6 6 3 2 16 4 6 3 -6 11 14 7 19 5 5 17 12 7 5 8 4 4 10 5 17 23 3 26 16 11 6 -2 6 8 4 11 17 19 8 12 15 6 12 56 2"

# --- OFF sheet: updated season totals ---
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("B2").Value = 12
$offWs.Range("C2").Value = 488
$offWs.Range("E2").Value = 14
$offWs.Range("F2").Value = 151
$offWs.Range("G2").Value = 118
$offWs.Range("H2").Value = 8
$offWs.Range("J2").Value = 63
$offWs.Range("L2").Value = 529
$offWs.Range("M2").Value = 353
$offWs.Range("Q2").Value = 1044
$offWs.Range("B3").Value = 19
$offWs.Range("C3").Value = 254
$offWs.Range("D3").Value = 10
$offWs.Range("E3").Value = 63
$offWs.Range("F3").Value = 189
$offWs.Range("H3").Value = 53
$offWs.Range("I3").Value = 106
$offWs.Range("J3").Value = 114
$offWs.Range("N3").Value = 36

# --- DEF sheet: updated season totals ---
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 343
$defWs.Range("D2").Value = 19
$defWs.Range("F2").Value = 110
$defWs.Range("G2").Value = 121
$defWs.Range("J2").Value = 63
$defWs.Range("L2").Value = 644
$defWs.Range("M2").Value = 430
$defWs.Range("O2").Value = 45
$defWs.Range("P2").Value = 27
$defWs.Range("Q2").Value = 1101
$defWs.Range("B3").Value = 24
$defWs.Range("C3").Value = 381
$defWs.Range("E3").Value = 57
$defWs.Range("F3").Value = 223
$defWs.Range("G3").Value = 62
$defWs.Range("H3").Value = 44
$defWs.Range("I3").Value = 124
$defWs.Range("J3").Value = 111
$defWs.Range("N3").Value = 32

# --- ST sheet: updated KO/PT totals and RM/D/RA logs ---
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 174
$stWs.Range("D2").Value = 103
$stWs.Range("F2").Value = 262
$stWs.Range("G2").Value = 249
$stWs.Range("J2").Value = 123
$stWs.Range("K2").Value = 116
$stWs.Range("L2").Value = 81
$stWs.Range("M2").Value = 63
$stWs.Range("B3").Value = 101
$stWs.Range("D3").Value = "41 66 50 47 34 31 44 59 40 55 32 41 54 63 59 61 59 36 36 55 42 60 48 65 47 46 41 47 17 41 52 47 48 32 33 54 42 42 37 50 40 38 47 32 41 41 60 46 51 53 50 49 25 47 59 45 41 53 42 12 63 30 42 41 39 45 44 42 63 44 45 49 50 52 37 49 41 57 35 42 42 36 41 42 46 40 27 53 47 48 35 37 35 39 46 50 39 54 53 39 46 33 59"
$stWs.Range("D4").Value = "0 12 0 7 0 0 0 0 0 -1 0 0 0 11 57 15 0 0 0 8 10 0 9 19 6 12 0 0 0 0 7 0 0 0 0 14 1 15 0 0 0 0 -1 0 0 0 0 2 0 0 0 0 0 11 13 0 0 0 0 0 14 0 18 0 0 0 0 0 6 14 0 0 11 16 0 0 0 0 0 0 7 0 4 0 0 0 0 10 15 8 0 0 0 0 0 17 11 11 0 0 12 0 0"
$stWs.Range("D5").Value = "11 9 -1 0 0 0 0 3 0 2 40 8 0 16 4 0 0 5 0 3 0 16 21 12 0 0 0 2 0 0 0 11 0 12 17 12 0 0 0 0 6 0 13 0 0 6 0 0 2 0 0 0 0 0 0 10 13 9 12 18 0 14 16 22 0 0 0 6 0 0 2 0 0 0 0 2 19 0 0 11 0 0 0 0 9 0 0 0 16 0 4 0 0 0 2 0 0 1 10 0 55 0 0 0 0 0"
$stWs.Range("B6").Value = "19 23 2 17 12 30 13 24 24 23 18 23 17 11 19 23 14 35 25 14 14 19 18 27 14 19 23 28 17 26 28 31 23 20 17 16 25 21 14 19 14 19 13 19 16 15 16 12 17 17 20 17 21 31 17 17 25 25 27 20 7 14 26 17 22 20 22 30 16 25 21 19 37 20"

# --- TURNS sheet: updated turnover counts ---
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("C2").Value = 17
$turnsWs.Range("E2").Value = 17
$turnsWs.Range("D3").Value = 12
$turnsWs.Range("E3").Value = 15

# --- PEN sheet: updated penalty counts ---
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value = 32
$penWs.Range("B3").Value = 27

